$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated run for publication: replace the 0/1 frequency matrix (rows 2-5,
# columns B-X) with the recomputed fractional frequencies.

$values = @(
  @(0.00337837837837838, 0, 0.981418918918919, 0.0118243243243243, 0, 0.00337837837837838, 0, 0, 0.0126689189189189, 0.997466216216216, 0.990709459459459, 0.00337837837837838, 0.000844594594594595, 0, 0.0177364864864865, 0.0152027027027027, 0.995777027027027, 0.980574324324324, 0.00591216216216216, 0.994932432432432, 0.983952702702703, 0.993243243243243, 0.0135135135135135),
  @(0.983952702702703, 0.980574324324324, 0, 0.983108108108108, 0.0135135135135135, 0, 0.995777027027027, 0.00168918918918919, 0.000844594594594595, 0.00168918918918919, 0.00253378378378378, 0.00506756756756757, 0.993243243243243, 0.981418918918919, 0.00337837837837838, 0.97972972972973, 0, 0.00168918918918919, 0.00168918918918919, 0.00253378378378378, 0.00337837837837838, 0.00506756756756757, 0.00168918918918919),
  @(0.0126689189189189, 0.0135135135135135, 0.0135135135135135, 0.00168918918918919, 0, 0.996621621621622, 0.00422297297297297, 0.0118243243243243, 0.00253378378378378, 0.000844594594594595, 0.00506756756756757, 0.989864864864865, 0.00506756756756757, 0.0126689189189189, 0.00253378378378378, 0.00506756756756757, 0.00337837837837838, 0.0143581081081081, 0, 0.000844594594594595, 0.00760135135135135, 0.000844594594594595, 0.97972972972973),
  @(0, 0.00506756756756757, 0.00506756756756757, 0.00337837837837838, 0.986486486486487, 0, 0, 0.986486486486487, 0.983952702702703, 0, 0.00168918918918919, 0.00168918918918919, 0.000844594594594595, 0.00591216216216216, 0.976351351351351, 0, 0.000844594594594595, 0.00337837837837838, 0.992398648648649, 0.00168918918918919, 0.00506756756756757, 0.000844594594594595, 0.00506756756756757)
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        # Column B is index 2
        $colIndex = 2 + $j
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}

$wb.Save()
